$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.451.51"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.676.54"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5310"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2695"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06410"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.16%  "
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.694.49"
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.507"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5576"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "0.0₅8351"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "26.486.15"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.743"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.343"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1289"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.397"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.442"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06333"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.638"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.453"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.679"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.010"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6182"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.14%  "
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.781"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01633"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "1.086.42"
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8646"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "1.821.29"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.163"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05206"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.481"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.042"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.11%  "
